# Generate Report for handoff
# Updates the localization-status workbook to reflect a failed handoff
# transform: new source file name, new status text, and removal of the
# (now non-existent) handoff target file / datetime info.

$wb = $excel.ActiveWorkbook

$newFileName = "6dd68bf3-7360-44d0-b101-438cfcb5d96a.md"
$newStatus   = "Handoff transform failed"
$naDate      = "0001-01-01 00:00:00"
$ignored     = "Ignored"

foreach ($name in @("Overview", "zh-cn", "de-de")) {
    $ws = $wb.Worksheets.Item($name)

    # Column A / B on row 2: new source file name + new status
    $ws.Range("A2").Value2 = $newFileName
    $ws.Range("B2").Value2 = $newStatus

    # Update the hyperlink display text that goes with A2 (keep the same
    # target address; only the visible file name text changes)
    foreach ($hl in $ws.Hyperlinks) {
        if ($hl.Range.Address() -eq '$A$2') {
            $hl.TextToDisplay = $newFileName
        }
    }
}

# The Overview sheet mirrors the status value in column C as well
# (same shared string as column B), so it must follow along too.
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("C2").Value2 = $newStatus

foreach ($name in @("zh-cn", "de-de")) {
    $ws = $wb.Worksheets.Item($name)

    # The handoff target file (C2) and its hyperlink no longer exist,
    # since the transform failed - remove the cell entirely.
    $ws.Range("C2").Clear()

    $toDelete = @()
    foreach ($hl in $ws.Hyperlinks) {
        if ($hl.Range.Address() -eq '$C$2') {
            $toDelete += $hl
        }
    }
    foreach ($hl in $toDelete) {
        $hl.Delete()
    }

    # Latest Target File datetime (D2) reverts to the "never happened"
    # sentinel date, same as row 3 below it.
    $ws.Range("D2").Value2 = $naDate

    # Handoff Reason (G2) / Dependency (H2) reset to defaults.
    $ws.Range("G2").Value2 = $naDate
    $ws.Range("H2").Value2 = $ignored

    # Row 3 Handoff Reason / Dependency columns also normalize to the
    # same shared values used elsewhere.
    $ws.Range("G3").Value2 = $naDate
    $ws.Range("H3").Value2 = $ignored
}
